$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.689.87"
$ws.Range("E2").Value = "  +2.80%  "

$ws.Range("D3").Value = "1.864.99"
$ws.Range("E3").Value = "  +2.10%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'246.50"
$ws.Range("E5").Value = "  +2.81%  "

$ws.Range("D6").Value = "'0.7014"
$ws.Range("E6").Value = "  +2.16%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'0.07776"
$ws.Range("E8").Value = "  +2.15%  "

$ws.Range("E9").Value = "  +2.16%  "

$ws.Range("E10").Value = "  +1.45%  "

$ws.Range("D11").Value = "'0.07832"
$ws.Range("E11").Value = "  +1.32%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.867.39"
$ws.Range("E12").Value = "  +2.23%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.179"
$ws.Range("E13").Value = "  +2.83%  "

$ws.Range("D14").Value = "'92.95"
$ws.Range("E14").Value = "  +3.34%  "

$ws.Range("D15").Value = "'0.6962"
$ws.Range("E15").Value = "  +3.76%  "

$ws.Range("D16").Value = "'6.637"
$ws.Range("E16").Value = "  +3.74%  "

$ws.Range("D17").Value = "29.687.75"
$ws.Range("E17").Value = "  +2.83%  "

$ws.Range("D18").Value = "'0.000008397"
$ws.Range("E18").Value = "  +1.51%  "

$ws.Range("D19").Value = "2.115.76"
$ws.Range("E19").Value = "  +1.94%  "

$ws.Range("D20").Value = "'243.74"
$ws.Range("E20").Value = "  +0.49%  "

$ws.Range("E21").Value = "  +1.70%  "

$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("D23").Value = "'7.659"
$ws.Range("E23").Value = "  +3.74%  "

$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("D25").Value = "'0.1523"
$ws.Range("E25").Value = "  +3.93%  "

$ws.Range("D26").Value = "'8.973"
$ws.Range("E26").Value = "  +3.22%  "

$ws.Range("D27").Value = "'160.20"
$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("D28").Value = "'18.39"
$ws.Range("E28").Value = "  +1.60%  "

$ws.Range("D29").Value = "'1.547"
$ws.Range("E29").Value = "  +1.07%  "

$ws.Range("E30").Value = "  +2.10%  "

$ws.Range("D31").Value = "'4.213"
$ws.Range("E31").Value = "  +1.63%  "

$ws.Range("D32").Value = "'1.201"
$ws.Range("E32").Value = "  +1.17%  "

$ws.Range("E33").Value = "  +0.48%  "

$ws.Range("D34").Value = "'0.7892"
$ws.Range("E34").Value = "  +4.78%  "

$ws.Range("D35").Value = "'1.918"
$ws.Range("E35").Value = "  +6.04%  "

$ws.Range("D36").Value = "'1.161"
$ws.Range("E36").Value = "  +2.01%  "

$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("D38").Value = "1.341.95"
$ws.Range("E38").Value = "  +11.01%  "

$ws.Range("D39").Value = "'0.01886"
$ws.Range("E39").Value = "  +3.07%  "

$ws.Range("D40").Value = "'2.744"
$ws.Range("E40").Value = "  +2.60%  "

$ws.Range("D41").Value = "'0.9615"
$ws.Range("E41").Value = "  +5.63%  "

$ws.Range("D42").Value = "'6.041"
$ws.Range("E42").Value = "  +13.57%  "

$ws.Range("D43").Value = "'106.73"
$ws.Range("E43").Value = "  -1.71%  "

$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("E45").Value = "  +3.68%  "

$ws.Range("D46").Value = "'9.792"
$ws.Range("E46").Value = "  +4.11%  "

$ws.Range("D47").Value = "2.012.48"
$ws.Range("E47").Value = "  +1.22%  "

$ws.Range("D48").Value = "'65.34"
$ws.Range("E48").Value = "  +4.66%  "

$ws.Range("D49").Value = "'0.5211"
$ws.Range("E49").Value = "  +0.89%  "

$ws.Range("D50").Value = "'1.795"
$ws.Range("E50").Value = "  +4.37%  "

$ws.Range("D51").Value = "'7.034"
$ws.Range("E51").Value = "  +2.89%  "
